# Update Betfair back/lay odds values on Sheet1 to reflect refreshed feed data.
# Only numeric cell values change; no rows/columns are inserted or removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 10).Value = 1.02
$ws.Cells.Item(2, 12).Value = 1.01
$ws.Cells.Item(2, 13).Value = 1.01
$ws.Cells.Item(2, 14).Value = 1.25
$ws.Cells.Item(2, 15).Value = 1.01
$ws.Cells.Item(2, 16).Value = 1.25
$ws.Cells.Item(2, 18).Value = 1.15
$ws.Cells.Item(2, 19).Value = 1.36
$ws.Cells.Item(2, 20).Value = 1.01
$ws.Cells.Item(2, 21).Value = 1.01
$ws.Cells.Item(2, 22).Value = 1.01
$ws.Cells.Item(2, 23).Value = 1.01
$ws.Cells.Item(2, 24).Value = 1000
$ws.Cells.Item(2, 25).Value = 1000
$ws.Cells.Item(2, 26).Value = 1000
$ws.Cells.Item(2, 27).Value = 1000
$ws.Cells.Item(2, 28).Value = 1000
$ws.Cells.Item(2, 29).Value = 1000
$ws.Cells.Item(2, 30).Value = 1000
$ws.Cells.Item(2, 31).Value = 1000
$ws.Cells.Item(2, 32).Value = 1000
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 34).Value = 1000
$ws.Cells.Item(2, 35).Value = 1000
$ws.Cells.Item(2, 36).Value = 1000
$ws.Cells.Item(2, 37).Value = 1000
$ws.Cells.Item(2, 38).Value = 1000
$ws.Cells.Item(2, 39).Value = 1000
$ws.Cells.Item(2, 40).Value = 1000
$ws.Cells.Item(2, 41).Value = 1000

# Row 3
$ws.Cells.Item(3, 16).Value = 1.74
$ws.Cells.Item(3, 17).Value = 1.93
$ws.Cells.Item(3, 22).Value = 1.58

# Row 4
$ws.Cells.Item(4, 6).Value = 1.66
$ws.Cells.Item(4, 7).Value = 1.97
$ws.Cells.Item(4, 8).Value = 2.04
$ws.Cells.Item(4, 9).Value = 14.5
$ws.Cells.Item(4, 10).Value = 3.1

# Row 6
$ws.Cells.Item(6, 6).Value = 1.21
$ws.Cells.Item(6, 7).Value = 1.28
$ws.Cells.Item(6, 8).Value = 15.5
$ws.Cells.Item(6, 9).Value = 23
$ws.Cells.Item(6, 10).Value = 6.6
$ws.Cells.Item(6, 11).Value = 9.199999999999999
$ws.Cells.Item(6, 16).Value = 2.54
$ws.Cells.Item(6, 17).Value = 1.51

# Row 7
$ws.Cells.Item(7, 17).Value = 1.89

# Row 10
$ws.Cells.Item(10, 17).Value = 1.97

# Row 11
$ws.Cells.Item(11, 7).Value = 1.33
$ws.Cells.Item(11, 16).Value = 3.05

# Row 16
$ws.Cells.Item(16, 8).Value = 8.4
$ws.Cells.Item(16, 9).Value = 9
$ws.Cells.Item(16, 20).Value = 2.34
$ws.Cells.Item(16, 21).Value = 1.7

# Row 17
$ws.Cells.Item(17, 14).Value = 2.28
$ws.Cells.Item(17, 15).Value = 1.53
$ws.Cells.Item(17, 19).Value = 5.3
$ws.Cells.Item(17, 20).Value = 2.04
$ws.Cells.Item(17, 24).Value = 9

# Row 18
$ws.Cells.Item(18, 7).Value = 2.08

# Row 20
$ws.Cells.Item(20, 6).Value = 2.14
$ws.Cells.Item(20, 7).Value = 2.3
$ws.Cells.Item(20, 8).Value = 3.45
$ws.Cells.Item(20, 9).Value = 3.8
$ws.Cells.Item(20, 10).Value = 3.45
$ws.Cells.Item(20, 17).Value = 1.89
$ws.Cells.Item(20, 21).Value = 2.18
$ws.Cells.Item(20, 22).Value = 1.36
$ws.Cells.Item(20, 23).Value = 1.76
$ws.Cells.Item(20, 25).Value = 15
$ws.Cells.Item(20, 26).Value = 28
$ws.Cells.Item(20, 27).Value = 75
$ws.Cells.Item(20, 28).Value = 11
$ws.Cells.Item(20, 29).Value = 9.800000000000001
$ws.Cells.Item(20, 30).Value = 16
$ws.Cells.Item(20, 31).Value = 48
$ws.Cells.Item(20, 32).Value = 15.5
$ws.Cells.Item(20, 33).Value = 12
$ws.Cells.Item(20, 34).Value = 18
$ws.Cells.Item(20, 36).Value = 980
$ws.Cells.Item(20, 37).Value = 25
$ws.Cells.Item(20, 38).Value = 980
$ws.Cells.Item(20, 40).Value = 18
$ws.Cells.Item(20, 41).Value = 48
